$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 151-152; existing rows 151-247 shift down to 153-249
$ws.Range("A151:A152").EntireRow.Insert()

# New row 151 (Primera)
$ws.Range("A151").Value = 1
$ws.Range("B151").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C151").Value = "Arica y Parinacota"
$ws.Range("D151").Value = 44596
$ws.Range("E151").Value = 15
$ws.Range("F151").Value = 100112043
$ws.Range("G151").Value = "Pepino ensalada"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 130
$ws.Range("K151").Value = 9000
$ws.Range("L151").Value = 10000
$ws.Range("M151").Value = 9500
$ws.Range("N151").Value = "`$/caja 70 unidades"
$ws.Range("O151").Value = "Región de Arica y Parinacota"
$ws.Range("P151").Value = 136
$ws.Range("Q151").Value = 70
$ws.Range("R151").Value = "Hortaliza"

# New row 152 (Segunda)
$ws.Range("A152").Value = 1
$ws.Range("B152").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C152").Value = "Arica y Parinacota"
$ws.Range("D152").Value = 44596
$ws.Range("E152").Value = 15
$ws.Range("F152").Value = 100112043
$ws.Range("G152").Value = "Pepino ensalada"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Segunda"
$ws.Range("J152").Value = 130
$ws.Range("K152").Value = 7000
$ws.Range("L152").Value = 8000
$ws.Range("M152").Value = 7500
$ws.Range("N152").Value = "`$/caja 100 unidades"
$ws.Range("O152").Value = "Región de Arica y Parinacota"
$ws.Range("P152").Value = 75
$ws.Range("Q152").Value = 100
$ws.Range("R152").Value = "Hortaliza"
